$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 1617488.4  # H17 was 1617492.8
$ws.Cells.Item(17, 10).Value = 1617488.4  # J17 was 1617492.8
$ws.Cells.Item(17, 12).Value = 4852465.199999999  # L17 was 4852478.4
$ws.Cells.Item(17, 14).Value = -4852801.199999999  # N17 was -4852814.4
$ws.Cells.Item(70, 8).Value = 2450  # H70 was 2362.2104
$ws.Cells.Item(70, 9).Value = 942.8570999999999  # I70 was 932.13336
$ws.Cells.Item(70, 11).Value = 2828.5713  # K70 was 2796.40008
$ws.Cells.Item(70, 13).Value = -2558.5713  # M70 was -2526.40008
$ws.Cells.Item(73, 8).Value = 2450  # H73 was 2362.2104
$ws.Cells.Item(73, 9).Value = 942.8570999999999  # I73 was 932.13336
$ws.Cells.Item(73, 11).Value = 2828.5713  # K73 was 2796.40008
$ws.Cells.Item(73, 13).Value = -1892.5713  # M73 was -1860.40008
$ws.Cells.Item(76, 8).Value = 3204.4  # H76 was 3208.5715
$ws.Cells.Item(76, 9).Value = 2819.0476  # I76 was 2820
$ws.Cells.Item(76, 10).Value = 5227.5  # J76 was 4935.5557
$ws.Cells.Item(76, 11).Value = 2819.0476  # K76 was 2820
$ws.Cells.Item(76, 12).Value = 5227.5  # L76 was 4935.5557
$ws.Cells.Item(76, 13).Value = -2504.0476  # M76 was -2505
$ws.Cells.Item(76, 14).Value = -5857.5  # N76 was -5565.5557
$ws.Cells.Item(79, 8).Value = 3204.4  # H79 was 3208.5715
$ws.Cells.Item(79, 9).Value = 2819.0476  # I79 was 2820
$ws.Cells.Item(79, 10).Value = 5227.5  # J79 was 4935.5557
$ws.Cells.Item(79, 11).Value = 2819.0476  # K79 was 2820
$ws.Cells.Item(79, 12).Value = 5227.5  # L79 was 4935.5557
$ws.Cells.Item(79, 13).Value = -1727.0476  # M79 was -1728
$ws.Cells.Item(79, 14).Value = -7411.5  # N79 was -7119.5557
$ws.Cells.Item(87, 8).Value = 15343.793  # H87 was 15245.179
$ws.Cells.Item(87, 10).Value = 15343.793  # J87 was 15245.179
$ws.Cells.Item(87, 12).Value = 15343.793  # L87 was 15245.179
$ws.Cells.Item(87, 14).Value = -17839.793  # N87 was -17741.179
$ws.Cells.Item(90, 8).Value = 15343.793  # H90 was 15245.179
$ws.Cells.Item(90, 10).Value = 15343.793  # J90 was 15245.179
$ws.Cells.Item(90, 12).Value = 46031.379  # L90 was 45735.537
$ws.Cells.Item(90, 14).Value = -58511.379  # N90 was -58215.537
$ws.Cells.Item(96, 8).Value = 20834796  # H96 was 15626221
$ws.Cells.Item(96, 9).Value = 41667936  # I96 was 31251010
$ws.Cells.Item(96, 10).Value = 1654.8334  # J96 was 1432.25
$ws.Cells.Item(96, 11).Value = 125003808  # K96 was 93753030
$ws.Cells.Item(96, 12).Value = 4964.5002  # L96 was 4296.75
$ws.Cells.Item(96, 13).Value = -125002435  # M96 was -93751657
$ws.Cells.Item(96, 14).Value = -7710.5002  # N96 was -7042.75
$ws.Cells.Item(112, 8).Value = 1223.4138  # H112 was 1241.4445
$ws.Cells.Item(112, 10).Value = 1283.96  # J112 was 1310.3914
$ws.Cells.Item(112, 12).Value = 3851.88  # L112 was 3931.1742
$ws.Cells.Item(112, 14).Value = -6067.88  # N112 was -6147.174199999999
$ws.Cells.Item(138, 8).Value = 4489.63  # H138 was 4187.5
$ws.Cells.Item(138, 9).Value = 2993.8064  # I138 was 2808.4285
$ws.Cells.Item(138, 10).Value = 5161.6665  # J138 was 4930.077
$ws.Cells.Item(138, 11).Value = 8981.4192  # K138 was 8425.2855
$ws.Cells.Item(138, 12).Value = 15484.9995  # L138 was 14790.231
$ws.Cells.Item(138, 13).Value = -3841.4192  # M138 was -3285.2855
$ws.Cells.Item(138, 14).Value = -25764.9995  # N138 was -25070.231

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(97, 8).Value = 717  # H97 was 759.16
$ws.Cells.Item(97, 9).Value = 447.7143  # I97 was 462
$ws.Cells.Item(97, 10).Value = 1007  # J97 was 1137.3636
$ws.Cells.Item(97, 11).Value = 447.7143  # K97 was 462
$ws.Cells.Item(97, 12).Value = 1007  # L97 was 1137.3636
$ws.Cells.Item(97, 13).Value = 48.28570000000002  # M97 was 34
$ws.Cells.Item(97, 14).Value = -1999  # N97 was -2129.3636
$ws.Cells.Item(124, 8).Value = 37559.8  # H124 was 22952.889
$ws.Cells.Item(124, 10).Value = 37559.8  # J124 was 22952.889
$ws.Cells.Item(124, 12).Value = 37559.8  # L124 was 22952.889
$ws.Cells.Item(124, 14).Value = -47379.8  # N124 was -32772.889
$ws.Cells.Item(139, 8).Value = 56232.5  # H139 was 60178.89
$ws.Cells.Item(139, 10).Value = 56232.5  # J139 was 60178.89
$ws.Cells.Item(139, 12).Value = 56232.5  # L139 was 60178.89
$ws.Cells.Item(139, 14).Value = -66512.5  # N139 was -70458.89

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(36, 8).Value = 1145.875  # H36 was 1000
$ws.Cells.Item(36, 9).Value = 565.2  # I36 was 1000
$ws.Cells.Item(36, 10).Value = 2113.6667  # J36 was 0
$ws.Cells.Item(36, 11).Value = 565.2  # K36 was 1000
$ws.Cells.Item(36, 12).Value = 2113.6667  # L36 was 0
$ws.Cells.Item(36, 13).Value = -31.20000000000005  # M36 was -466
$ws.Cells.Item(36, 14).Value = -3181.6667  # N36 was None
$ws.Cells.Item(86, 8).Value = 3486.6667  # H86 was 2797.1365
$ws.Cells.Item(86, 9).Value = 4016.6667  # I86 was 2726.0908
$ws.Cells.Item(86, 10).Value = 3133.3333  # J86 was 2868.182
$ws.Cells.Item(86, 11).Value = 4016.6667  # K86 was 2726.0908
$ws.Cells.Item(86, 12).Value = 3133.3333  # L86 was 2868.182
$ws.Cells.Item(86, 13).Value = -2893.6667  # M86 was -1603.0908
$ws.Cells.Item(86, 14).Value = -5379.3333  # N86 was -5114.182
$ws.Cells.Item(89, 8).Value = 3486.6667  # H89 was 2797.1365
$ws.Cells.Item(89, 9).Value = 4016.6667  # I89 was 2726.0908
$ws.Cells.Item(89, 10).Value = 3133.3333  # J89 was 2868.182
$ws.Cells.Item(89, 11).Value = 20083.3335  # K89 was 13630.454
$ws.Cells.Item(89, 12).Value = 15666.6665  # L89 was 14340.91
$ws.Cells.Item(89, 13).Value = -14467.3335  # M89 was -8014.454
$ws.Cells.Item(89, 14).Value = -26898.6665  # N89 was -25572.91
$ws.Cells.Item(140, 8).Value = 59185  # H140 was 59653.332
$ws.Cells.Item(140, 10).Value = 59185  # J140 was 59653.332
$ws.Cells.Item(140, 12).Value = 59185  # L140 was 59653.332
$ws.Cells.Item(140, 14).Value = -69545  # N140 was -70013.33199999999

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(138, 8).Value = 40297.668  # H138 was 40336.555
$ws.Cells.Item(138, 10).Value = 40297.668  # J138 was 40336.555
$ws.Cells.Item(138, 12).Value = 40297.668  # L138 was 40336.555
$ws.Cells.Item(138, 14).Value = -50577.668  # N138 was -50616.555
$ws.Cells.Item(140, 8).Value = 74298.44500000001  # H140 was 72257.125
$ws.Cells.Item(140, 10).Value = 74298.44500000001  # J140 was 72257.125
$ws.Cells.Item(140, 12).Value = 74298.44500000001  # L140 was 72257.125
$ws.Cells.Item(140, 14).Value = -84658.44500000001  # N140 was -82617.125

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(54, 8).Value = 4133.3335  # H54 was 0
$ws.Cells.Item(54, 10).Value = 4133.3335  # J54 was 0
$ws.Cells.Item(54, 12).Value = 12400.0005  # L54 was 0
$ws.Cells.Item(54, 14).Value = -13518.0005  # N54 was None

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 2450  # H80 was 2900
$ws.Cells.Item(80, 9).Value = 2000  # I80 was 0
$ws.Cells.Item(80, 11).Value = 2000  # K80 was 0
$ws.Cells.Item(80, 13).Value = -1002  # M80 was None
$ws.Cells.Item(83, 8).Value = 2450  # H83 was 2900
$ws.Cells.Item(83, 9).Value = 2000  # I83 was 0
$ws.Cells.Item(83, 11).Value = 10000  # K83 was 0
$ws.Cells.Item(83, 13).Value = -5008  # M83 was None
$ws.Cells.Item(97, 8).Value = 2766.182  # H97 was 2295.8333
$ws.Cells.Item(97, 9).Value = 2936.5557  # I97 was 2295.8333
$ws.Cells.Item(97, 10).Value = 1999.5  # J97 was 0
$ws.Cells.Item(97, 11).Value = 2936.5557  # K97 was 2295.8333
$ws.Cells.Item(97, 12).Value = 1999.5  # L97 was 0
$ws.Cells.Item(97, 13).Value = -2440.5557  # M97 was -1799.8333
$ws.Cells.Item(97, 14).Value = -2991.5  # N97 was None
$ws.Cells.Item(122, 8).Value = 73639.66  # H122 was 2423.2415
$ws.Cells.Item(122, 9).Value = 105806.71  # I122 was 1725.9
$ws.Cells.Item(122, 10).Value = 3457  # J122 was 3972.889
$ws.Cells.Item(122, 11).Value = 317420.13  # K122 was 5177.700000000001
$ws.Cells.Item(122, 12).Value = 10371  # L122 was 11918.667
$ws.Cells.Item(122, 13).Value = -314970.13  # M122 was -2727.700000000001
$ws.Cells.Item(122, 14).Value = -15271  # N122 was -16818.667
$ws.Cells.Item(135, 8).Value = 0  # H135 was 40600
$ws.Cells.Item(135, 10).Value = 0  # J135 was 40600
$ws.Cells.Item(135, 12).Value = 0  # L135 was 40600
$ws.Cells.Item(135, 14).ClearContents()  # N135 removed (was -50740)
$ws.Cells.Item(138, 8).Value = 49955.3  # H138 was 48884.93
$ws.Cells.Item(138, 10).Value = 49955.3  # J138 was 48884.93
$ws.Cells.Item(138, 12).Value = 49955.3  # L138 was 48884.93
$ws.Cells.Item(138, 14).Value = -60235.3  # N138 was -59164.93

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 1293.8235  # H22 was 1363
$ws.Cells.Item(22, 9).Value = 914.5  # I22 was 949.375
$ws.Cells.Item(22, 11).Value = 914.5  # K22 was 949.375
$ws.Cells.Item(22, 13).Value = -619.5  # M22 was -654.375
$ws.Cells.Item(27, 8).Value = 1293.8235  # H27 was 1363
$ws.Cells.Item(27, 9).Value = 914.5  # I27 was 949.375
$ws.Cells.Item(27, 11).Value = 914.5  # K27 was 949.375
$ws.Cells.Item(27, 13).Value = -807.5  # M27 was -842.375
$ws.Cells.Item(46, 8).Value = 1743.4706  # H46 was 1680.0555
$ws.Cells.Item(46, 10).Value = 1780.5  # J46 was 1649.5555
$ws.Cells.Item(46, 12).Value = 1780.5  # L46 was 1649.5555
$ws.Cells.Item(46, 14).Value = -2156.5  # N46 was -2025.5555
$ws.Cells.Item(68, 8).Value = 1775.7778  # H68 was 1896
$ws.Cells.Item(68, 9).Value = 1550.5  # I68 was 1500
$ws.Cells.Item(68, 10).Value = 1956  # J68 was 1995
$ws.Cells.Item(68, 11).Value = 1550.5  # K68 was 1500
$ws.Cells.Item(68, 12).Value = 1956  # L68 was 1995
$ws.Cells.Item(68, 13).Value = -801.5  # M68 was -751
$ws.Cells.Item(68, 14).Value = -3454  # N68 was -3493
$ws.Cells.Item(71, 8).Value = 1775.7778  # H71 was 1896
$ws.Cells.Item(71, 9).Value = 1550.5  # I71 was 1500
$ws.Cells.Item(71, 10).Value = 1956  # J71 was 1995
$ws.Cells.Item(71, 11).Value = 7752.5  # K71 was 1500
$ws.Cells.Item(71, 12).Value = 9780  # L71 was 9975
$ws.Cells.Item(71, 13).Value = -4008.5  # M71 was -3756
$ws.Cells.Item(71, 14).Value = -17268  # N71 was -17463
$ws.Cells.Item(82, 8).Value = 2672.1667  # H82 was 2635.375
$ws.Cells.Item(82, 9).Value = 3745  # I82 was 2830
$ws.Cells.Item(82, 10).Value = 2135.75  # J82 was 2518.6
$ws.Cells.Item(82, 11).Value = 3745  # K82 was 2830
$ws.Cells.Item(82, 12).Value = 2135.75  # L82 was 2518.6
$ws.Cells.Item(82, 13).Value = -3384  # M82 was -2469
$ws.Cells.Item(82, 14).Value = -2857.75  # N82 was -3240.6
$ws.Cells.Item(85, 8).Value = 2672.1667  # H85 was 2635.375
$ws.Cells.Item(85, 9).Value = 3745  # I85 was 2830
$ws.Cells.Item(85, 10).Value = 2135.75  # J85 was 2518.6
$ws.Cells.Item(85, 11).Value = 3745  # K85 was 2830
$ws.Cells.Item(85, 12).Value = 2135.75  # L85 was 2518.6
$ws.Cells.Item(85, 13).Value = -2497  # M85 was -1582
$ws.Cells.Item(85, 14).Value = -4631.75  # N85 was -5014.6
$ws.Cells.Item(100, 8).Value = 1633.3334  # H100 was 3860
$ws.Cells.Item(100, 9).Value = 1866.6666  # I100 was 2850
$ws.Cells.Item(100, 10).Value = 1400  # J100 was 4533.3335
$ws.Cells.Item(100, 11).Value = 1866.6666  # K100 was 2850
$ws.Cells.Item(100, 12).Value = 1400  # L100 was 4533.3335
$ws.Cells.Item(100, 13).Value = -1325.6666  # M100 was -2309
$ws.Cells.Item(100, 14).Value = -2482  # N100 was -5615.3335
$ws.Cells.Item(122, 8).Value = 37038596  # H122 was 15875950
$ws.Cells.Item(122, 9).Value = 111111110  # I122 was 55556896
$ws.Cells.Item(122, 10).Value = 2340  # J122 was 3572
$ws.Cells.Item(122, 11).Value = 333333330  # K122 was 166670688
$ws.Cells.Item(122, 12).Value = 7020  # L122 was 10716
$ws.Cells.Item(122, 13).Value = -333330880  # M122 was -166668238
$ws.Cells.Item(122, 14).Value = -11920  # N122 was -15616
$ws.Cells.Item(133, 8).Value = 27751  # H133 was 23260
$ws.Cells.Item(133, 10).Value = 27751  # J133 was 23260
$ws.Cells.Item(133, 12).Value = 27751  # L133 was 23260
$ws.Cells.Item(133, 14).Value = -32811  # N133 was -28320

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 113731.11  # H122 was 127722.5
$ws.Cells.Item(122, 9).Value = 145368.58  # I122 was 169296.67
$ws.Cells.Item(122, 11).Value = 436105.74  # K122 was 507890.01
$ws.Cells.Item(122, 13).Value = -433655.74  # M122 was -505440.01
$ws.Cells.Item(137, 8).Value = 55131.4  # H137 was 56050.4
$ws.Cells.Item(137, 10).Value = 55131.4  # J137 was 56050.4
$ws.Cells.Item(137, 12).Value = 55131.4  # L137 was 56050.4
$ws.Cells.Item(137, 14).Value = -65331.4  # N137 was -66250.39999999999
$ws.Cells.Item(140, 8).Value = 51666.332  # H140 was 35000
$ws.Cells.Item(140, 9).Value = 45000  # I140 was 25000
$ws.Cells.Item(140, 10).Value = 54999.5  # J140 was 55000
$ws.Cells.Item(140, 11).Value = 45000  # K140 was 25000
$ws.Cells.Item(140, 12).Value = 54999.5  # L140 was 55000
$ws.Cells.Item(140, 13).Value = -39820  # M140 was -19820
$ws.Cells.Item(140, 14).Value = -65359.5  # N140 was -65360
$ws.Cells.Item(141, 8).Value = 130857.5  # H141 was 78507.5
$ws.Cells.Item(141, 9).Value = 0  # I141 was 20000
$ws.Cells.Item(141, 10).Value = 130857.5  # J141 was 90209
$ws.Cells.Item(141, 11).Value = 0  # K141 was 20000
$ws.Cells.Item(141, 12).Value = 130857.5  # L141 was 90209
$ws.Cells.Item(141, 13).ClearContents()  # M141 removed (was -14820)
$ws.Cells.Item(141, 14).Value = -141217.5  # N141 was -100569
